$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi4"
$ws.Range("C2").Value = "Adam23"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03034166666666667
$ws.Range("H2").Value = 0.091025
$ws.Range("I2").Value = 0.002265018822860075
$ws.Range("J2").Value = 0.002265018822860075
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03926266666666667
$ws.Range("N2").Value = 0.117788
$ws.Range("O2").Value = 0.005313231574131687
$ws.Range("P2").Value = 0.005313231574131686
$ws.Range("Q2").Value = 0.001191294744444445
$ws.Range("R2").Value = 0.0107216527
$ws.Range("S2").Value = 0.00001203456952562274
$ws.Range("T2").Value = 0.00001203456952562274

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi4"
$ws.Range("C3").Value = "Adam23"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03034166666666667
$ws.Range("H3").Value = 0.091025
$ws.Range("I3").Value = 0.002265018822860075
$ws.Range("J3").Value = 0.002265018822860075
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.402094666666667
$ws.Range("N3").Value = 13.206284
$ws.Range("O3").Value = 0.5957147173375057
$ws.Range("P3").Value = 0.5957147173375056
$ws.Range("Q3").Value = 0.1335668890111111
$ws.Range("R3").Value = 1.2021020011
$ws.Range("S3").Value = 0.00134930504782422
$ws.Range("T3").Value = 0.001349305047824219

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi4"
$ws.Range("C4").Value = "Adam23"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03034166666666667
$ws.Range("H4").Value = 0.091025
$ws.Range("I4").Value = 0.002265018822860075
$ws.Range("J4").Value = 0.002265018822860075
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.948244666666667
$ws.Range("N4").Value = 8.844734
$ws.Range("O4").Value = 0.3989720510883627
$ws.Range("P4").Value = 0.3989720510883626
$ws.Range("Q4").Value = 0.08945465692777779
$ws.Range("R4").Value = 0.80509191235
$ws.Range("S4").Value = 0.0009036792055102332
$ws.Range("T4").Value = 0.0009036792055102329

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Lgi4"
$ws.Range("C5").Value = "Adam23"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.326286
$ws.Range("H5").Value = 30.978858
$ws.Range("I5").Value = 0.7708618124768957
$ws.Range("J5").Value = 0.7708618124768957
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03926266666666667
$ws.Range("N5").Value = 0.117788
$ws.Range("O5").Value = 0.005313231574131687
$ws.Range("P5").Value = 0.005313231574131686
$ws.Range("Q5").Value = 0.4054375251226667
$ws.Range("R5").Value = 3.648937726104
$ws.Range("S5").Value = 0.004095767321344622
$ws.Range("T5").Value = 0.004095767321344621

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Lgi4"
$ws.Range("C6").Value = "Adam23"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.326286
$ws.Range("H6").Value = 30.978858
$ws.Range("I6").Value = 0.7708618124768957
$ws.Range("J6").Value = 0.7708618124768957
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.402094666666667
$ws.Range("N6").Value = 13.206284
$ws.Range("O6").Value = 0.5957147173375057
$ws.Range("P6").Value = 0.5957147173375056
$ws.Range("Q6").Value = 45.45728852707467
$ws.Range("R6").Value = 409.115596743672
$ws.Range("S6").Value = 0.4592137267259512
$ws.Range("T6").Value = 0.4592137267259511

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Lgi4"
$ws.Range("C7").Value = "Adam23"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.326286
$ws.Range("H7").Value = 30.978858
$ws.Range("I7").Value = 0.7708618124768957
$ws.Range("J7").Value = 0.7708618124768957
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.948244666666667
$ws.Range("N7").Value = 8.844734
$ws.Range("O7").Value = 0.3989720510883627
$ws.Range("P7").Value = 0.3989720510883626
$ws.Range("Q7").Value = 30.44441762597467
$ws.Range("R7").Value = 273.999758633772
$ws.Range("S7").Value = 0.3075523184295999
$ws.Range("T7").Value = 0.3075523184295998

# Row 8
$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Lgi4"
$ws.Range("C8").Value = "Adam23"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.039140333333334
$ws.Range("H8").Value = 9.117421
$ws.Range("I8").Value = 0.2268731687002442
$ws.Range("J8").Value = 0.2268731687002443
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.03926266666666667
$ws.Range("N8").Value = 0.117788
$ws.Range("O8").Value = 0.005313231574131687
$ws.Range("P8").Value = 0.005313231574131686
$ws.Range("Q8").Value = 0.1193247538608889
$ws.Range("R8").Value = 1.073922784748
$ws.Range("S8").Value = 0.001205429683261442
$ws.Range("T8").Value = 0.001205429683261442

# Row 9
$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Lgi4"
$ws.Range("C9").Value = "Adam23"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.039140333333334
$ws.Range("H9").Value = 9.117421
$ws.Range("I9").Value = 0.2268731687002442
$ws.Range("J9").Value = 0.2268731687002443
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.402094666666667
$ws.Range("N9").Value = 13.206284
$ws.Range("O9").Value = 0.5957147173375057
$ws.Range("P9").Value = 0.5957147173375056
$ws.Range("Q9").Value = 13.37858345261822
$ws.Range("R9").Value = 120.407251073564
$ws.Range("S9").Value = 0.1351516855637303
$ws.Range("T9").Value = 0.1351516855637302

# Row 10
$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Lgi4"
$ws.Range("C10").Value = "Adam23"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.039140333333334
$ws.Range("H10").Value = 9.117421
$ws.Range("I10").Value = 0.2268731687002442
$ws.Range("J10").Value = 0.2268731687002443
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.948244666666667
$ws.Range("N10").Value = 8.844734
$ws.Range("O10").Value = 0.3989720510883627
$ws.Range("P10").Value = 0.3989720510883626
$ws.Range("Q10").Value = 8.960129279001558
$ws.Range("R10").Value = 80.64116351101401
$ws.Range("S10").Value = 0.09051605345325257
$ws.Range("T10").Value = 0.09051605345325256
